$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewUser")

# Update the test-data row (row 2) with the new registration/reset-password
# test credentials: a fresh email, password and EID.
$ws.Range("A2").Value = "TestPF1221+13082021092544@gmail.com"
$ws.Range("B2").Value = "pfqa_123"
$ws.Range("C2").Value = "TestPF1221_13082021092544"
